$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Degree")

# Update header labels for B1 (DegreeAbrev) and C1 (DegreeName)
$ws.Range("B1").Value = "DegreeAbrev"
$ws.Range("C1").Value = "DegreeName "

# Fill D3:D5 and E3:E5 with the same opening/closing literals already used in D2/E2
$ws.Range("D3").Value = "new Degree{"
$ws.Range("E3").Value = "},"
$ws.Range("D4").Value = "new Degree{"
$ws.Range("E4").Value = "},"
$ws.Range("D5").Value = "new Degree{"
$ws.Range("E5").Value = "},"

# Update/add the concatenation formula for F2:F5 (now also includes DegreeAbrev/DegreeName)
$ws.Range("F2").Formula = '=D2&$A$1&"="&A2&","&$B$1&"="&B2&","&$C$1&"="&C2&E2'
$ws.Range("F3").Formula = '=D3&$A$1&"="&A3&","&$B$1&"="&B3&","&$C$1&"="&C3&E3'
$ws.Range("F4").Formula = '=D4&$A$1&"="&A4&","&$B$1&"="&B4&","&$C$1&"="&C4&E4'
$ws.Range("F5").Formula = '=D5&$A$1&"="&A5&","&$B$1&"="&B5&","&$C$1&"="&C5&E5'

# Resize columns D and F to fit the new content
$ws.Range("D1").ColumnWidth = 13
$ws.Range("F1").ColumnWidth = 67.5

# Move the active selection as in the edited workbook
$ws.Range("F15").Select()
